$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Province"
$ws.Range("B1").Value = "Region"
$ws.Range("C1").Value = "Police Service"
$ws.Range("D1").Value = "Source"
$ws.Range("E1").Value = "Format"
$ws.Range("F1").Value = "Open"
$ws.Range("G1").Value = "Special Note"

# --- Row 2 (QC) ---
$ws.Range("A2").Value = "QC"
$ws.Range("B2").Value = "Adminstrative Regions"
$ws.Range("C2").Value = "Surêté Québec"
$ws.Range("D2").Value = "https://mern.gouv.qc.ca/territoire/portrait/portrait-donnees-mille.jsp"
$ws.Range("E2").Value = "Shapefile"
$ws.Range("F2").Value = "Open"

# --- Row 3 (NS) ---
$ws.Range("A3").Value = "NS"
$ws.Range("B3").Value = "Counties"
$ws.Range("C3").Value = "???"
$ws.Range("D3").Value = "https://data.novascotia.ca/browse?q=Nova%20Scotia%20Topographic%20Database%20-%20County%20Boundaries&sortBy=relevance"
$ws.Range("E3").Value = "Various"
$ws.Range("F3").Value = "Open"

# --- Row 4 (ON) ---
$ws.Range("A4").Value = "ON"
$ws.Range("B4").Value = "Counties"
$ws.Range("C4").Value = "Ontario Provincial Police"
$ws.Range("D4").Value = "https://www.ontario.ca/data/municipal-boundaries"
$ws.Range("E4").Value = "Shapefile"
$ws.Range("F4").Value = "Open"

# --- Row 5 (CAN) ---
$ws.Range("A5").Value = "CAN"
$ws.Range("B5").Value = "Detachments"
$ws.Range("C5").Value = "RCMP"
$ws.Range("D5").Value = "http://www.rcmp-grc.gc.ca/detach/en/d/697/kml"
$ws.Range("E5").Value = "KML"
$ws.Range("F5").Value = "???"
$ws.Range("G5").Value = "Needs scraping"

# --- Column widths (target "best fit" widths: 21.42578125, 23.140625, 122.5703125,
#     11.140625, 17.5703125, 20.140625 characters respectively) ---
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws.Columns.Item(3).ColumnWidth = 22.333333333333332
$ws.Columns.Item(4).ColumnWidth = 121.66666666666667
$ws.Columns.Item(5).ColumnWidth = 10.333333333333334
$ws.Columns.Item(6).ColumnWidth = 16.666666666666668
$ws.Columns.Item(7).ColumnWidth = 19.333333333333332

# --- Header row bottom border (thin bottom border under row 1) ---
$headerRange = $ws.Range("A1:G1")
$headerRange.Borders.Item(9).LineStyle = 1

# --- Selection ---
$ws.Range("G5").Select()
